# Update cryptocurrency price (D) and volume-change (E) cells per latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.935.19"
$ws.Range("E2").Value = "  -2.49%  "
$ws.Range("D3").Value = "2.097.91"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  -1.15%  "
$ws.Range("D5").Value = "'346.28"
$ws.Range("E5").Value = "  +2.29%  "
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("D7").Value = "'0.5147"
$ws.Range("E7").Value = "  -2.22%  "
$ws.Range("D8").Value = "'0.4408"
$ws.Range("E8").Value = "  -3.37%  "
$ws.Range("D9").Value = "'0.09366"
$ws.Range("E9").Value = "  +2.56%  "
$ws.Range("D10").Value = "'52.03"
$ws.Range("E10").Value = "  -5.46%  "
$ws.Range("D11").Value = "'1.167"
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D13").Value = "2.096.39"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").Value = "'6.727"
$ws.Range("E14").Value = "  -2.09%  "
$ws.Range("D15").Value = "'8.143"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "'99.20"
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("D17").Value = "'0.00001161"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").Value = "'20.62"
$ws.Range("E19").Value = "  +5.25%  "
$ws.Range("D20").Value = "'0.06672"
$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("D22").Value = "'6.207"
$ws.Range("E22").Value = "  -2.12%  "
$ws.Range("D23").Value = "30.033.88"
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("E24").Value = "  -2.82%  "
$ws.Range("D25").Value = "'2.330"
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("D26").Value = "2.351.38"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").Value = "'21.94"
$ws.Range("E27").Value = "  -2.66%  "
$ws.Range("D28").Value = "'2.552"
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("D29").Value = "'162.32"
$ws.Range("E29").Value = "  -2.05%  "
$ws.Range("D30").Value = "'133.02"
$ws.Range("E30").Value = "  -1.93%  "
$ws.Range("D31").Value = "'1.167"
$ws.Range("D32").Value = "'0.1059"
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("D33").Value = "'1.635"
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("D34").Value = "'6.210"
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("D35").Value = "'3.954"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").Value = "'6.201"
$ws.Range("E36").Value = "  +4.92%  "
$ws.Range("D37").Value = "'10.14"
$ws.Range("E37").Value = "  -4.76%  "
$ws.Range("E38").Value = "  -4.08%  "
$ws.Range("D39").Value = "'0.06782"
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("E40").Value = "  -2.79%  "
$ws.Range("D41").Value = "'12.48"
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("D42").Value = "'0.6901"
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("D43").Value = "'1.301"
$ws.Range("E43").Value = "  +3.16%  "
$ws.Range("D44").Value = "'0.6605"
$ws.Range("E44").Value = "  +1.40%  "
$ws.Range("D45").Value = "'14.10"
$ws.Range("E45").Value = "  -7.53%  "
$ws.Range("D46").Value = "'2.272"
$ws.Range("E46").Value = "  -2.18%  "
$ws.Range("D47").Value = "'3.632"
$ws.Range("E47").Value = "  -1.94%  "
$ws.Range("D48").Value = "'0.00000000355"
$ws.Range("E48").Value = "  -5.89%  "
$ws.Range("D49").Value = "'1.216"
$ws.Range("E49").Value = "  -3.32%  "
$ws.Range("D50").Value = "'81.74"
$ws.Range("E50").Value = "  -2.46%  "
$ws.Range("D51").Value = "'0.07193"
$ws.Range("E51").Value = "  -1.75%  "
